$d = $word.ActiveDocument

# 1. Merge the split runs' text into a single contiguous run of text by
#    replacing the old fragmented phrase with the same text (Find/Replace
#    collapses formatting-identical adjacent runs into one run).
$d.Content.Find.Execute(
    "Incidencias (y soluci" + [char]0x00F3 + "n en el caso de ser necesario) en la construcci" + [char]0x00F3 + "n del componente HW del",
    $true, $false, $false, $false, $false,
    $true, 1, $false,
    "Incidencias (y soluci" + [char]0x00F3 + "n en el caso de ser necesario) en la construcci" + [char]0x00F3 + "n del componente HW del",
    2
)

# 2. Set the paragraph mark's run properties (pPr/rPr) for the first
#    paragraph: Arial font, size 25 (12.5pt), single underline.
$para = $d.Paragraphs.Item(1)
$markRange = $para.Range
$markRange.SetRange($markRange.End - 1, $markRange.End)
$markRange.Font.Name = "Arial"
$markRange.Font.Size = 12.5
$markRange.Font.Underline = 1
